$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an exact literal string into a cell without Excel's
# "looks like a number" auto-conversion, by round-tripping through a
# quote-prefixed (text) scratch cell via Copy / PasteSpecial (values).
function Set-ExactText($cellAddr, $text) {
    $ws.Range("ZZ1").Value = "'" + $text
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-ExactText "D2" "69.224.18"
$ws.Range("E2").Value = "  -2.12%  "
Set-ExactText "D3" "3.485.20"
$ws.Range("E3").Value = "  -2.23%  "
Set-ExactText "D4" "1.00"
$ws.Range("E4").Value = "  -0.11%  "
Set-ExactText "D5" "611.65"
$ws.Range("E5").Value = "  +4.72%  "
Set-ExactText "D6" "185.47"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.76%  "
Set-ExactText "D8" "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -2.56%  "
Set-ExactText "D10" "0.649"
$ws.Range("E10").Value = "  -0.63%  "
Set-ExactText "D11" "53.02"
$ws.Range("E11").Value = "  -2.71%  "
Set-ExactText "D12" "0.0000306"
$ws.Range("E12").Value = "  -3.41%  "
Set-ExactText "D13" "9.52"
$ws.Range("E13").Value = "  +0.29%  "
Set-ExactText "D14" "4.049.94"
$ws.Range("E14").Value = "  -1.97%  "
Set-ExactText "D15" "601.12"
$ws.Range("E15").Value = "  +4.49%  "
Set-ExactText "D16" "69.405.07"
$ws.Range("E16").Value = "  -1.96%  "
Set-ExactText "D17" "12.59"
$ws.Range("E17").Value = "  +1.30%  "
Set-ExactText "D18" "18.79"
$ws.Range("E18").Value = "  -2.59%  "
Set-ExactText "D19" "3.493.95"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("E20").Value = "  -0.40%  "
Set-ExactText "D21" "0.986"
$ws.Range("E21").Value = "  -1.76%  "
Set-ExactText "D22" "17.23"
$ws.Range("E22").Value = "  -2.75%  "
Set-ExactText "D23" "104.89"
$ws.Range("E23").Value = "  +9.79%  "
Set-ExactText "D24" "4.65"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +0.07%  "
Set-ExactText "D26" "3.01"
$ws.Range("E26").Value = "  +2.29%  "
Set-ExactText "D27" "10.92"
$ws.Range("E27").Value = "  -3.00%  "
Set-ExactText "D28" "9.92"
$ws.Range("E28").Value = "  +8.54%  "
Set-ExactText "D29" "33.45"
$ws.Range("E29").Value = "  +3.01%  "
Set-ExactText "D30" "6.96"
$ws.Range("E30").Value = "  -3.81%  "
Set-ExactText "D31" "12.38"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  -0.17%  "
Set-ExactText "D33" "3.91"
$ws.Range("E33").Value = "  +16.47%  "
Set-ExactText "D34" "63.31"
$ws.Range("E34").Value = "  +0.32%  "
Set-ExactText "D35" "3.16"
$ws.Range("E35").Value = "  -6.89%  "
Set-ExactText "D36" "0.998"
$ws.Range("E36").Value = "  -0.18%  "
Set-ExactText "D37" "520.10"
$ws.Range("E37").Value = "  -5.32%  "
Set-ExactText "D38" "0.396"
$ws.Range("E38").Value = "  -4.60%  "
Set-ExactText "D39" "3.575.88"
$ws.Range("E39").Value = "  +0.23%  "
Set-ExactText "D40" "3.57"
$ws.Range("E40").Value = "  +4.50%  "
Set-ExactText "D41" "36.60"
$ws.Range("E41").Value = "  -3.35%  "
Set-ExactText "D42" "0.0₃0772"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  +0.72%  "
Set-ExactText "D44" "0.0459"
$ws.Range("E44").Value = "  +2.88%  "
Set-ExactText "D45" "2.95"
$ws.Range("E45").Value = "  +0.70%  "
Set-ExactText "D46" "0.143"
$ws.Range("E46").Value = "  +3.45%  "
Set-ExactText "D47" "3.31"
$ws.Range("E47").Value = "  -6.35%  "
Set-ExactText "D48" "8.82"
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("E49").Value = "  +0.38%  "

# Coin list update: FLOKI newly inserted at rank 49 (row 50), pushing
# Monero down to row 51 with a refreshed price/volume; OceanProtocol
# (formerly row 51) drops off the bottom of the list entirely.
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-ExactText "D50" "0.000243"
$ws.Range("E50").Value = "  -8.19%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-ExactText "D51" "131.11"
$ws.Range("E51").Value = "  -2.57%  "

# Clean up the scratch column used by Set-ExactText.
$ws.Columns("ZZ:ZZ").Delete()
